$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing data (D:K) right to (F:M)
$ws.Range("D:E").Insert()

# Copy the (now shifted) number/date formats from F:G back onto the new D:E columns
# so the new columns look like their neighbours (date format on header rows, number format elsewhere)
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Approximate the original bestFit column widths for the two new columns
$ws.Columns("D:E").ColumnWidth = 13.8

# Set new values for columns D and E (two new most-recent quarters)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 2896000
$ws.Range("E8").Value = 2947000
$ws.Range("D9").Value = 1448000
$ws.Range("E9").Value = 1548000
$ws.Range("D10").Value = 1448000
$ws.Range("E10").Value = 1399000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 281000
$ws.Range("E15").Value = 276000
$ws.Range("D17").Value = 1818000
$ws.Range("E17").Value = 1927000
$ws.Range("D18").Value = 1078000
$ws.Range("E18").Value = 1020000
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 30000
$ws.Range("D21").Value = 1360000
$ws.Range("E21").Value = 1326000
$ws.Range("D22").Value = 148000
$ws.Range("E22").Value = 142000
$ws.Range("D23").Value = 930000
$ws.Range("E23").Value = 908000
$ws.Range("D24").Value = 228000
$ws.Range("E24").Value = 206000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 702000
$ws.Range("E26").Value = 702000
$ws.Range("D27").Value = 700000
$ws.Range("E27").Value = 700000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = -30000
$ws.Range("D33").Value = 700000
$ws.Range("E33").Value = 700000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 700000
$ws.Range("E35").Value = 700000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 358000
$ws.Range("E41").Value = 729000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 1009000
$ws.Range("E43").Value = 1043000
$ws.Range("D44").Value = 207000
$ws.Range("E44").Value = 267000
$ws.Range("D45").Value = 288000
$ws.Range("E45").Value = 70000
$ws.Range("D46").Value = 1862000
$ws.Range("E46").Value = 2109000
$ws.Range("D47").Value = 3109000
$ws.Range("E47").Value = 3109000
$ws.Range("D48").Value = 31091000
$ws.Range("E48").Value = 30712000
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 177000
$ws.Range("E52").Value = 392000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 36239000
$ws.Range("E54").Value = 36322000
$ws.Range("D57").Value = 1505000
$ws.Range("E57").Value = 1394000
$ws.Range("D58").Value = 585000
$ws.Range("E58").Value = 500000
$ws.Range("D59").Value = 501000
$ws.Range("E59").Value = 547000
$ws.Range("D60").Value = 2591000
$ws.Range("E60").Value = 2441000
$ws.Range("D61").Value = 10560000
$ws.Range("E61").Value = 10635000
$ws.Range("D62").Value = 7726000
$ws.Range("E62").Value = 7766000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 20877000
$ws.Range("E66").Value = 20842000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 13440000
$ws.Range("E72").Value = 13645000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 15362000
$ws.Range("E76").Value = 15480000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 700000
$ws.Range("E81").Value = 700000
$ws.Range("D83").Value = 282000
$ws.Range("E83").Value = 276000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 842000
$ws.Range("E89").Value = 1058000
$ws.Range("D91").Value = -625000
$ws.Range("E91").Value = -490000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -517000
$ws.Range("E94").Value = -355000
$ws.Range("D96").Value = -217000
$ws.Range("E96").Value = -219000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -696000
$ws.Range("E100").Value = -316000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -371000
$ws.Range("E102").Value = 387000

Write-Host "done"
